$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table of events and
# need the same updates applied to column F (想去人数 / "want to go" count).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1584
    $ws.Range("F3").Value = 208
    $ws.Range("F4").Value = 195
    $ws.Range("F5").Value = 3385
    $ws.Range("F6").Value = 5952
    $ws.Range("F7").Value = 322
    $ws.Range("F8").Value = 30
    $ws.Range("F11").Value = 8789
    $ws.Range("F12").Value = 2353
    $ws.Range("F13").Value = 247
    $ws.Range("F14").Value = 5349
    $ws.Range("F15").Value = 10248
}
